# Update countries & provincias Spain
#
# - Reorders "Guinea-Bisau" in the country table so it appears right after
#   "Maldivas" (before "Hong Kong"), and refreshes its stats.
# - Refreshes the "Datos actualizados ..." timestamp caption.
# - Updates case counters for several countries (Estados Unidos, row 14,
#   Hong Kong/Tunez/Guinea-Bisau block, Sri Lanka).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- A1: refresh "last updated" caption -------------------------------
$ws.Range("A1").Value = "Datos actualizados a 20 de Mayo de 2020 a las 21:35"

# --- Row 4: Estados Unidos --------------------------------------------
$ws.Range("B4").Value = 1579985
$ws.Range("C4").Value = 9402
$ws.Range("D4").Value = 365747
$ws.Range("E4").Value = 1120057
$ws.Range("G4").Value = 648
$ws.Range("H4").Value = 94181

# --- Row 14 -------------------------------------------------------------
$ws.Range("B14").Value = 112015
$ws.Range("C14").Value = 5540
$ws.Range("E14").Value = 63159

# --- Rows 102-104: move Guinea-Bisau ahead of Hong Kong / Tunez --------
# Row 102 becomes Guinea-Bisau with its updated stats.
$ws.Range("A102").Value = "Guinea-Bisau"
$ws.Range("B102").Value = 1089
$ws.Range("C102").Value = 51
$ws.Range("D102").Value = 42
$ws.Range("E102").Value = 1041
$ws.Range("F102").Value = 0
$ws.Range("G102").Value = 0
$ws.Range("H102").Value = 6

# Row 103 becomes Hong Kong (its stats are unchanged, just shifted down).
$ws.Range("A103").Value = "Hong Kong"
$ws.Range("B103").Value = 1056
$ws.Range("C103").Value = 0
$ws.Range("D103").Value = 1026
$ws.Range("E103").Value = 26
$ws.Range("F103").Value = 0
$ws.Range("G103").Value = 0
$ws.Range("H103").Value = 4

# Row 104 becomes Tunez (its stats are unchanged, just shifted down).
$ws.Range("A104").Value = "Tunez"
$ws.Range("B104").Value = 1044
$ws.Range("C104").Value = 1
$ws.Range("D104").Value = 826
$ws.Range("E104").Value = 171
$ws.Range("F104").Value = 0
$ws.Range("G104").Value = 1
$ws.Range("H104").Value = 47

# Row 105 (Kenia) is unchanged.

# --- Row 106: Sri Lanka --------------------------------------------------
$ws.Range("B106").Value = 1028
$ws.Range("C106").Value = 5
$ws.Range("E106").Value = 435
